# "Generate Report for Handoff"
# A new handoff batch was generated: the priority of 4 files changed from
# "low" to "ht", and the "Latest Handoff Datetime" batch timestamps for
# those same 4 files were refreshed, on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column (E) for rows 4-7 changes from "low" to "ht" on both sheets
$wsZhCn.Range("E4:E7").Value = "ht"
$wsDeDe.Range("E4:E7").Value = "ht"

# Latest Handoff Datetime column (H) for rows 4-7 is refreshed with the new
# handoff-generation timestamp for each locale
$wsZhCn.Range("H4:H7").Value = "2016-09-06 04:39:35"
$wsDeDe.Range("H4:H7").Value = "2016-09-06 04:39:41"
